# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the latest data refresh.
# Price values are forced to text (leading apostrophe) so strings like "27.096.81"
# or "1.0000" are preserved exactly as text instead of being parsed as numbers,
# then the cell style is reset to Normal so no stray quote-prefix style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.Value = "'27.096.81"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "

$dCell = $ws.Range("D3")
$dCell.Value = "'1.889.88"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +1.41%  "

$dCell = $ws.Range("D4")
$dCell.Value = "'1.001"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$dCell = $ws.Range("D5")
$dCell.Value = "'306.73"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

$dCell = $ws.Range("D6")
$dCell.Value = "'1.0000"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$dCell = $ws.Range("D7")
$dCell.Value = "'0.5143"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "

$dCell = $ws.Range("D8")
$dCell.Value = "'0.3757"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +3.16%  "

$dCell = $ws.Range("D9")
$dCell.Value = "'0.07204"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +0.32%  "

$dCell = $ws.Range("D10")
$dCell.Value = "'21.17"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "

$dCell = $ws.Range("D11")
$dCell.Value = "'0.9022"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "

$dCell = $ws.Range("D12")
$dCell.Value = "'0.07653"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +2.13%  "

$dCell = $ws.Range("D13")
$dCell.Value = "'1.874.27"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "

$dCell = $ws.Range("D14")
$dCell.Value = "'94.48"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "

$dCell = $ws.Range("D15")
$dCell.Value = "'5.250"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "

$dCell = $ws.Range("D16")
$dCell.Value = "'1.001"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "

$dCell = $ws.Range("D17")
$dCell.Value = "'0.000008487"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "

$dCell = $ws.Range("D18")
$dCell.Value = "'14.43"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +1.66%  "

$dCell = $ws.Range("D19")
$dCell.Value = "'1.000"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "

$dCell = $ws.Range("D20")
$dCell.Value = "'27.123.15"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

$dCell = $ws.Range("D21")
$dCell.Value = "'5.061"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "

$dCell = $ws.Range("D22")
$dCell.Value = "'2.125.27"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +1.14%  "

$dCell = $ws.Range("D23")
$dCell.Value = "'10.57"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +1.88%  "

$dCell = $ws.Range("D24")
$dCell.Value = "'6.390"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "

$dCell = $ws.Range("D25")
$dCell.Value = "'2.297"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +11.14%  "

$dCell = $ws.Range("D26")
$dCell.Value = "'147.38"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

$dCell = $ws.Range("D27")
$dCell.Value = "'1.767"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "

$dCell = $ws.Range("D28")
$dCell.Value = "'18.05"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "

$dCell = $ws.Range("D29")
$dCell.Value = "'114.14"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "

$dCell = $ws.Range("D30")
$dCell.Value = "'4.936"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +5.33%  "

$dCell = $ws.Range("D31")
$dCell.Value = "'4.795"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +2.37%  "

$dCell = $ws.Range("D32")
$dCell.Value = "'0.09196"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -0.70%  "

$dCell = $ws.Range("D33")
$dCell.Value = "'0.05068"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -0.68%  "

$dCell = $ws.Range("D34")
$dCell.Value = "'1.237"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +7.42%  "

$dCell = $ws.Range("D35")
$dCell.Value = "'0.7723"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +2.79%  "

$dCell = $ws.Range("D36")
$dCell.Value = "'2.993"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "

$dCell = $ws.Range("D37")
$dCell.Value = "'3.281"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

$dCell = $ws.Range("D38")
$dCell.Value = "'2.610"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  +2.38%  "

$dCell = $ws.Range("D39")
$dCell.Value = "'0.5590"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  +0.93%  "

$dCell = $ws.Range("D40")
$dCell.Value = "'0.01992"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -0.56%  "

$dCell = $ws.Range("D41")
$dCell.Value = "'1.074"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "

$dCell = $ws.Range("D42")
$dCell.Value = "'9.070"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +6.73%  "

$dCell = $ws.Range("D43")
$dCell.Value = "'6.644"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +1.70%  "

$dCell = $ws.Range("D44")
$dCell.Value = "'117.83"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

$dCell = $ws.Range("D45")
$dCell.Value = "'0.1503"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +2.15%  "

$dCell = $ws.Range("D46")
$dCell.Value = "'0.4804"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "

$dCell = $ws.Range("D47")
$dCell.Value = "'10.18"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +1.49%  "

$ws.Range("E48").Value = "  +0.11%  "

$dCell = $ws.Range("D49")
$dCell.Value = "'1.596"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +1.86%  "

$dCell = $ws.Range("D50")
$dCell.Value = "'37.57"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +2.01%  "

$dCell = $ws.Range("D51")
$dCell.Value = "'64.03"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
